$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.18588399887085
$ws.Range("B1").Value = 2.115474939346313
$ws.Range("C1").Value = 6.472019672393799
$ws.Range("D1").Value = 2.305414915084839
$ws.Range("E1").Value = 1.195563793182373
